$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "rule exchange" config table appended below the existing tables,
# starting at row 128 (B:D), mirroring the layout/styling of the other
# tables already on the sheet.
$data = @(
    @(1,  "ID",                               "NUMBER (19)"),
    @(2,  "CREATED_AT",                       "TIMESTAMP"),
    @(3,  "CREATED_BY",                       "VARCHAR2 (255 CHAR)"),
    @(4,  "IS_DELETED",                       "NUMBER (1)"),
    @(5,  "UPDATED_AT",                       "TIMESTAMP"),
    @(6,  "UPDATED_BY",                       "VARCHAR2 (255 CHAR)"),
    @(7,  "EXCHANGE_POINT",                   "NUMBER (19)"),
    @(8,  "EXCHANGE_VALUE",                   "NUMBER (19)"),
    @(9,  "FIX_POINT_AMOUNT",                 "NUMBER (19)"),
    @(10, "FREQUENCY_LIMIT_EVENT_PER_USER",   "VARCHAR2 (255 CHAR)"),
    @(11, "FREQUENCY_LIMIT_POINT_PER_USER",   "VARCHAR2 (255 CHAR)"),
    @(12, "FREQUENCY_TIME_WAIT",              "VARCHAR2 (255 CHAR)"),
    @(13, "IS_EXCHANGE_BY_VALUE",             "NUMBER (1)"),
    @(14, "IS_NET_VALUE",                     "NUMBER (1)"),
    @(15, "LIMIT_EVENT_PER_USER",             "NUMBER (19)"),
    @(16, "LIMIT_POINT_PER_TRANSACTION",      "NUMBER (19)"),
    @(17, "LIMIT_POINT_PER_USER",             "NUMBER (19)"),
    @(18, "MIN_TRANSACTION",                  "NUMBER (19)"),
    @(19, "RULE_ID",                          "NUMBER (19)"),
    @(20, "TIME_WAIT",                        "NUMBER (19)")
)

$startRow = 128
$row = $startRow
foreach ($d in $data) {
    $ws.Range("B" + $row).Value = $d[0]
    $ws.Range("C" + $row).Value = $d[1]
    $ws.Range("D" + $row).Value = $d[2]
    $row = $row + 1
}
$endRow = $row - 1

$tableRange = $ws.Range("B" + $startRow + ":D" + $endRow)
$tableRange.WrapText = $true
$tableRange.VerticalAlignment = -4108

$selRange = $ws.Range("B" + $startRow + ":D" + $endRow)
$selRange.Select()
